# Re-ran "resolve" and "classify+summarise" steps after changes to the
# mapping file. This re-generated the summary numbers on a few sheets -
# several counts collapsed to 0 (and their accompanying percentage column
# was dropped), and the "High Priority break-up" sheet now only reports
# the IUCN breakdown (the old Range row's figures were folded into a
# single IUCN row).

$wb = $excel.ActiveWorkbook

# ---- "Range Status" sheet: counts all zeroed out, percentage column removed ----
$ws2 = $wb.Worksheets.Item("Range Status")
for ($r = 2; $r -le 7; $r++) {
    $ws2.Cells.Item($r, 2).Value = 0      # column B
    $ws2.Cells.Item($r, 3).ClearContents() # column C (no longer populated)
}

# ---- "Species qualification" sheet: Range Analysis count zeroed ----
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Cells.Item(5, 2).Value = 0   # B5, "Range Analysis" row

# ---- "High Priority break-up" sheet: now a single IUCN row ----
$ws5 = $wb.Worksheets.Item("High Priority break-up")
$ws5.Cells.Item(2, 1).Value = "IUCN"
$ws5.Cells.Item(2, 2).Value = 18
$ws5.Cells.Item(2, 3).Value = 100
$ws5.Cells.Item(2, 4).Value = 18
$ws5.Cells.Item(2, 5).Value = 100

# remove the old third row (previously the IUCN row, now unused)
$ws5.Rows.Item(3).ClearContents()
